$wb = $excel.ActiveWorkbook

# ALC row 31
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 300
$ws.Range("I31").Value = 300
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 900
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -670
$ws.Range("N31").Value = ""

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2156.1
$ws.Range("I33").Value = 1284.5555
$ws.Range("J33").Value = 10000
$ws.Range("K33").Value = 1284.5555
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = -1055.5555
$ws.Range("N33").Value = -10458

# ALC row 97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1900
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1900
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 5700
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = -6692

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2614.1428
$ws.Range("I113").Value = 2387.375
$ws.Range("J113").Value = 2916.5
$ws.Range("K113").Value = 2387.375
$ws.Range("L113").Value = 2916.5
$ws.Range("M113").Value = 866.625
$ws.Range("N113").Value = -9424.5

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1925.1765
$ws.Range("I131").Value = 917.5
$ws.Range("J131").Value = 3364.7144
$ws.Range("K131").Value = 2752.5
$ws.Range("L131").Value = 10094.1432
$ws.Range("M131").Value = 2287.5
$ws.Range("N131").Value = -20174.1432

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2337
$ws.Range("I138").Value = 1857.625
$ws.Range("K138").Value = 5572.875
$ws.Range("M138").Value = -432.875

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2808.625
$ws.Range("I88").Value = 2984.5
$ws.Range("J88").Value = 2750
$ws.Range("K88").Value = 2984.5
$ws.Range("L88").Value = 2750
$ws.Range("M88").Value = -2578.5
$ws.Range("N88").Value = -3562

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2808.625
$ws.Range("I91").Value = 2984.5
$ws.Range("J91").Value = 2750
$ws.Range("K91").Value = 2984.5
$ws.Range("L91").Value = 2750
$ws.Range("M91").Value = -1580.5
$ws.Range("N91").Value = -5558

# ARM row 101
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 75458.86
$ws.Range("J101").Value = 75458.86
$ws.Range("L101").Value = 75458.86
$ws.Range("N101").Value = -81948.86

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3642.182
$ws.Range("I132").Value = 2364.1538
$ws.Range("J132").Value = 5488.222
$ws.Range("K132").Value = 7092.4614
$ws.Range("L132").Value = 16464.666
$ws.Range("M132").Value = -4562.4614
$ws.Range("N132").Value = -21524.666

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 60715
$ws.Range("J139").Value = 60715
$ws.Range("L139").Value = 60715
$ws.Range("N139").Value = -70995

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2615.0715
$ws.Range("I86").Value = 2362.6
$ws.Range("J86").Value = 3246.25
$ws.Range("K86").Value = 2362.6
$ws.Range("L86").Value = 3246.25
$ws.Range("M86").Value = -1239.6
$ws.Range("N86").Value = -5492.25

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2615.0715
$ws.Range("I89").Value = 2362.6
$ws.Range("J89").Value = 3246.25
$ws.Range("K89").Value = 11813
$ws.Range("L89").Value = 16231.25
$ws.Range("M89").Value = -6197
$ws.Range("N89").Value = -27463.25

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4827.4854
$ws.Range("I31").Value = 1121.1842
$ws.Range("J31").Value = 9522.134
$ws.Range("K31").Value = 1121.1842
$ws.Range("L31").Value = 9522.134
$ws.Range("M31").Value = -826.1841999999999
$ws.Range("N31").Value = -10112.134

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4827.4854
$ws.Range("I34").Value = 1121.1842
$ws.Range("J34").Value = 9522.134
$ws.Range("K34").Value = 1121.1842
$ws.Range("L34").Value = 9522.134
$ws.Range("M34").Value = -919.1841999999999
$ws.Range("N34").Value = -9926.134

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8333.333000000001
$ws.Range("I62").Value = 7500
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 7500
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -6876
$ws.Range("N62").Value = -11248

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 8333.333000000001
$ws.Range("I65").Value = 7500
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 37500
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -34380
$ws.Range("N65").Value = -56240

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 24001950
$ws.Range("I4").Value = 17502062
$ws.Range("J4").Value = 50001500
$ws.Range("K4").Value = 52506186
$ws.Range("L4").Value = 150004500
$ws.Range("M4").Value = -52506074
$ws.Range("N4").Value = -150004724

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 732.5294
$ws.Range("I5").Value = 515.8125
$ws.Range("J5").Value = 4200
$ws.Range("K5").Value = 1547.4375
$ws.Range("L5").Value = 12600
$ws.Range("M5").Value = -1435.4375
$ws.Range("N5").Value = -12824

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 7969.5713
$ws.Range("I122").Value = 497.66666
$ws.Range("K122").Value = 4478.99994
$ws.Range("M122").Value = -2028.99994

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 732.5294
$ws.Range("I135").Value = 515.8125
$ws.Range("J135").Value = 4200
$ws.Range("K135").Value = 4642.3125
$ws.Range("L135").Value = 37800
$ws.Range("M135").Value = -2107.3125
$ws.Range("N135").Value = -42870

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5461.727
$ws.Range("I70").Value = 5365.0454
$ws.Range("J70").Value = 5655.091
$ws.Range("K70").Value = 5365.0454
$ws.Range("L70").Value = 5655.091
$ws.Range("M70").Value = -5095.0454
$ws.Range("N70").Value = -6195.091

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5461.727
$ws.Range("I73").Value = 5365.0454
$ws.Range("J73").Value = 5655.091
$ws.Range("K73").Value = 5365.0454
$ws.Range("L73").Value = 5655.091
$ws.Range("M73").Value = -4429.0454
$ws.Range("N73").Value = -7527.091

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3039
$ws.Range("I126").Value = 3085.3333
$ws.Range("K126").Value = 9255.999899999999
$ws.Range("M126").Value = -6785.999899999999

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2514.875
$ws.Range("I132").Value = 2033
$ws.Range("J132").Value = 3166.8235
$ws.Range("K132").Value = 6099
$ws.Range("L132").Value = 9500.470499999999
$ws.Range("M132").Value = -3569
$ws.Range("N132").Value = -14560.4705

# GSM row 137
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 45557.145
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 45557.145
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 45557.145
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -55757.145

# GSM row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 955
$ws.Range("I55").Value = 834
$ws.Range("J55").Value = 1076
$ws.Range("K55").Value = 834
$ws.Range("L55").Value = 1076
$ws.Range("M55").Value = -661
$ws.Range("N55").Value = -1422

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3218.5
$ws.Range("I132").Value = 2281.75
$ws.Range("J132").Value = 4155.25
$ws.Range("K132").Value = 6845.25
$ws.Range("L132").Value = 12465.75
$ws.Range("M132").Value = -4315.25
$ws.Range("N132").Value = -17525.75

# WVR row 43
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = 0

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1566.091
$ws.Range("I132").Value = 1293.8889
$ws.Range("J132").Value = 1998.4117
$ws.Range("K132").Value = 3881.6667
$ws.Range("L132").Value = 5995.2351
$ws.Range("M132").Value = -1351.6667
$ws.Range("N132").Value = -11055.2351
